$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.45
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 1.07
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 1.36
$ws.Range("M2").Value = 3.2
$ws.Range("N2").Value = 2.1
$ws.Range("O2").Value = 1.73
$ws.Range("R2").Value = 2.25
$ws.Range("S2").Value = 1.57
$ws.Range("U2").Value = 6.5
$ws.Range("W2").Value = 10
$ws.Range("AE2").Value = 15
$ws.Range("AG2").Value = 23
$ws.Range("AH2").Value = 101
$ws.Range("AI2").Value = 67
$ws.Range("AJ2").Value = 67

$ws.Range("N3").Value = 2.1
$ws.Range("O3").Value = 1.73

$ws.Range("G4").Value = 2.25
$ws.Range("I4").Value = 3.2
$ws.Range("N4").Value = 2.05
$ws.Range("O4").Value = 1.85
$ws.Range("T4").Value = 8.5
$ws.Range("Y4").Value = 29
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 41
$ws.Range("AJ4").Value = 41

$ws.Range("I5").Value = 8.5
$ws.Range("L5").Value = 1.11
$ws.Range("M5").Value = 6.5
$ws.Range("N5").Value = 1.4
$ws.Range("O5").Value = 2.88
$ws.Range("P5").Value = 1.22
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.95
$ws.Range("T5").Value = 10
$ws.Range("U5").Value = 8
$ws.Range("Y5").Value = 21
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 13
$ws.Range("AD5").Value = 201
$ws.Range("AE5").Value = 26
$ws.Range("AJ5").Value = 41

$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 3.9
$ws.Range("R6").Value = 1.67
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 8.5
$ws.Range("U6").Value = 10
$ws.Range("V6").Value = 8.5
$ws.Range("W6").Value = 17
$ws.Range("AE6").Value = 13
$ws.Range("AF6").Value = 21
$ws.Range("AG6").Value = 13
$ws.Range("AI6").Value = 29
$ws.Range("AJ6").Value = 34

$ws.Range("G7").Value = 1.62
$ws.Range("I7").Value = 5
$ws.Range("AA7").Value = 8
$ws.Range("AG7").Value = 17

$ws.Range("G8").Value = 2.7
$ws.Range("I8").Value = 2.8
$ws.Range("N8").Value = 2.88
$ws.Range("O8").Value = 1.4
$ws.Range("P8").Value = 1.67
$ws.Range("Q8").Value = 2.1
$ws.Range("T8").Value = 6
$ws.Range("U8").Value = 11

$ws.Range("G9").Value = 3.4
$ws.Range("I9").Value = 2.25
$ws.Range("AB9").Value = 21
$ws.Range("AF9").Value = 9

$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 4
$ws.Range("K10").Value = 12
$ws.Range("N10").Value = 1.8
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 1.36
$ws.Range("Q10").Value = 3
$ws.Range("U10").Value = 9.5
$ws.Range("X10").Value = 15
$ws.Range("AD10").Value = 151
$ws.Range("AF10").Value = 21
$ws.Range("AG10").Value = 13
$ws.Range("AI10").Value = 29

$ws.Range("N11").Value = 1.8
$ws.Range("O11").Value = 2

$ws.Range("G13").Value = 5.25
$ws.Range("H13").Value = 3.5
$ws.Range("I13").Value = 1.7
$ws.Range("J13").Value = 1.08
$ws.Range("K13").Value = 8
$ws.Range("R13").Value = 2.1
$ws.Range("S13").Value = 1.67
$ws.Range("U13").Value = 26
$ws.Range("Z13").Value = 8
$ws.Range("AA13").Value = 7
$ws.Range("AE13").Value = 5.5
$ws.Range("AF13").Value = 7
$ws.Range("AI13").Value = 15

$ws.Range("G15").Value = 2.35
$ws.Range("I15").Value = 3
$ws.Range("AH15").Value = 34

$ws.Range("N16").Value = 2.08
$ws.Range("O16").Value = 1.73

$ws.Range("N17").Value = 2.08
$ws.Range("O17").Value = 1.73

$ws.Range("K19").Value = 8.5

$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.98
$ws.Range("O20").Value = 1.88

$ws.Range("J21").Value = 1.03
$ws.Range("K21").Value = 15

$ws.Range("H22").Value = 3.7
$ws.Range("I22").Value = 2.25
$ws.Range("L22").Value = 1.25
$ws.Range("M22").Value = 3.75
$ws.Range("N22").Value = 1.88
$ws.Range("O22").Value = 1.98
$ws.Range("P22").Value = 1.36
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 1.67
$ws.Range("S22").Value = 2.1
$ws.Range("T22").Value = 10
$ws.Range("U22").Value = 15
$ws.Range("X22").Value = 23
$ws.Range("Y22").Value = 29
$ws.Range("Z22").Value = 12
$ws.Range("AA22").Value = 7
$ws.Range("AB22").Value = 15
$ws.Range("AD22").Value = 201
$ws.Range("AE22").Value = 8.5
$ws.Range("AF22").Value = 11
$ws.Range("AG22").Value = 9.5
$ws.Range("AJ22").Value = 26

$ws.Range("G23").Value = 2.75
$ws.Range("I23").Value = 2.35
$ws.Range("K23").Value = 12

$ws.Range("G24").Value = 3.5
$ws.Range("I24").Value = 2.25
$ws.Range("J24").Value = 1.11
$ws.Range("K24").Value = 6.5
$ws.Range("L24").Value = 1.5
$ws.Range("M24").Value = 2.5
$ws.Range("T24").Value = 8
$ws.Range("AD24").Value = 501
$ws.Range("AI24").Value = 21

$ws.Range("G25").Value = 1.9
$ws.Range("I25").Value = 3.9
$ws.Range("U25").Value = 9
$ws.Range("V25").Value = 8.5
$ws.Range("X25").Value = 15
$ws.Range("AF25").Value = 21
$ws.Range("AI25").Value = 34
$ws.Range("AJ25").Value = 41

$ws.Range("H27").Value = 3.5
$ws.Range("I27").Value = 2.75
$ws.Range("K27").Value = 15
$ws.Range("L27").Value = 1.18
$ws.Range("M27").Value = 4.5
$ws.Range("N27").Value = 1.65
$ws.Range("O27").Value = 2.2
$ws.Range("P27").Value = 1.33
$ws.Range("Q27").Value = 3.25
$ws.Range("R27").Value = 1.53
$ws.Range("S27").Value = 2.38
$ws.Range("T27").Value = 11
$ws.Range("U27").Value = 13
$ws.Range("V27").Value = 9.5
$ws.Range("Y27").Value = 23
$ws.Range("Z27").Value = 15
$ws.Range("AA27").Value = 7
$ws.Range("AB27").Value = 12
$ws.Range("AD27").Value = 126
$ws.Range("AE27").Value = 12
$ws.Range("AF27").Value = 15
$ws.Range("AI27").Value = 21

$ws.Range("I29").Value = 3.5
$ws.Range("J29").Value = 1.05
$ws.Range("K29").Value = 8
$ws.Range("R29").Value = 1.6
$ws.Range("T29").Value = 8.75
$ws.Range("U29").Value = 10.5
$ws.Range("V29").Value = 8.25
$ws.Range("W29").Value = 18
$ws.Range("Z29").Value = 8
$ws.Range("AA29").Value = 6.9
$ws.Range("AF29").Value = 21

Write-Output "Applied 194 cell updates"
